# Add Research Vocabularies Australia-era entry: ScoLOMFR
# (commit: "add ScoLOMFR")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 18 values, mirroring the style of the existing data rows (A/B as text, C/D as numeric 0 with style 2)
$ws.Range("A18").Value = "ScoLOMFR"
$ws.Range("B18").Value = "https://www.reseau-canope.fr/scolomfr/data"
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0

# Match the style used by the other data rows (C2:D17 use style index 2 i.e. the "Standard" cell style with fontId 0 applyFont)
$ws.Range("C18:D18").Style = $ws.Range("C17:D17").Style

# Update the active selection to the newly added cell, as captured in the sheetView
$ws.Range("B18").Select()

$wb.Save()
